$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths: target OOXML width 15.42578125 for both columns. The engine
# snaps ColumnWidth input to 1/6-character increments (pixel grid), so the
# closest achievable width is 15.5; 14.6 is an input that lands there.
$ws.Columns("A").ColumnWidth = 14.6
$ws.Columns("B").ColumnWidth = 14.6

$colA = @(
    -0.41111872142219852,
    -0.29047129344553468,
    -0.18394193417701388,
    -0.17078024383805968,
    -0.16372336043372826,
    -0.10947248384668473,
    -0.089285485139184217,
    -0.068813533873430366,
    -0.062405437455806734,
    -0.056057885437397204,
    -0.051506718304171528,
    -0.045422229629727084,
    -0.039161897089641684,
    -0.027090944963487296,
    -0.021056540394184253,
    -0.015029449014302143,
    -0.0090048199958552289,
    -0.058616618479181426,
    -0.049574178442955308,
    -0.018014040178234936,
    -0.0090043251229037224,
    -0.093955471220644071,
    -0.084639649455335864,
    -0.042127912980734372,
    -0.046121706634643544,
    -0.040051738224342159,
    -0.033967181461515228,
    -0.027696145469221278,
    -0.015521894358267119,
    0.0045468595348814311,
    -0.027020058273578229,
    -0.0060008691013333504
)
$colB = @(
    0.40981334584981255,
    0.28689740510936979,
    0.1827802437601882,
    0.16972336014835498,
    0.16159419018596211,
    0.10928548478845457,
    0.088813533518317556,
    0.068405437149253956,
    0.062057885126125534,
    0.056006717998151601,
    0.05142222931673146,
    0.045161896771564347,
    0.039090944620749468,
    0.027056540073330915,
    0.021029448692341468,
    0.015004819672430392,
    0.0089999996646294989,
    0.058574178163489421,
    0.049279265381536153,
    0.018004324838454266,
    0.0089999997152503397,
    0.093639649164542149,
    0.084127912558260576,
    0.04199999957517786,
    0.046051737928479497,
    0.039967181164914933,
    0.033696145170049707,
    0.02752189403445815,
    0.015453140110047059,
    -0.004570871055717074,
    0.027000868742218742,
    0.0059999996984299031
)

for ($i = 0; $i -lt 32; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $colA[$i]
    $ws.Cells.Item($i + 1, 2).Value = $colB[$i]
}
